# Update NATMI LR-pair output (Sema3f-Nrp1) with newly-computed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 6.704275
$ws.Range("H2").Value = 20.112825
$ws.Range("I2").Value = 0.4617710489234531
$ws.Range("J2").Value = 0.4617710489234532
$ws.Range("M2").Value = 123.2806423333333
$ws.Range("N2").Value = 369.841927
$ws.Range("O2").Value = 0.6241574062367528
$ws.Range("P2").Value = 0.6241574062367526
$ws.Range("Q2").Value = 826.5073283793083
$ws.Range("R2").Value = 7438.565955413776
$ws.Range("S2").Value = 0.2882178201712872
$ws.Range("T2").Value = 0.2882178201712871

$ws.Range("G3").Value = 6.704275
$ws.Range("H3").Value = 20.112825
$ws.Range("I3").Value = 0.4617710489234531
$ws.Range("J3").Value = 0.4617710489234532
$ws.Range("O3").Value = 0.2392728888301323
$ws.Range("P3").Value = 0.2392728888301322
$ws.Range("Q3").Value = 316.8444275827083
$ws.Range("R3").Value = 2851.599848244375
$ws.Range("S3").Value = 0.110489292854035
$ws.Range("T3").Value = 0.110489292854035

$ws.Range("G4").Value = 6.704275
$ws.Range("H4").Value = 20.112825
$ws.Range("I4").Value = 0.4617710489234531
$ws.Range("J4").Value = 0.4617710489234532
$ws.Range("O4").Value = 0.136569704933115
$ws.Range("P4").Value = 0.136569704933115
$ws.Range("Q4").Value = 180.8451855796833
$ws.Range("R4").Value = 1627.60667021715
$ws.Range("S4").Value = 0.06306393589813102
$ws.Range("T4").Value = 0.06306393589813102

$ws.Range("I5").Value = 0.03922895479591048
$ws.Range("J5").Value = 0.03922895479591048
$ws.Range("M5").Value = 123.2806423333333
$ws.Range("N5").Value = 369.841927
$ws.Range("O5").Value = 0.6241574062367528
$ws.Range("P5").Value = 0.6241574062367526
$ws.Range("Q5").Value = 70.21448984094999
$ws.Range("R5").Value = 631.93040856855
$ws.Range("S5").Value = 0.02448504267479431
$ws.Range("T5").Value = 0.0244850426747943

$ws.Range("I6").Value = 0.03922895479591048
$ws.Range("J6").Value = 0.03922895479591048
$ws.Range("O6").Value = 0.2392728888301323
$ws.Range("P6").Value = 0.2392728888301322
$ws.Range("S6").Value = 0.009386425339804171
$ws.Range("T6").Value = 0.009386425339804171

$ws.Range("I7").Value = 0.03922895479591048
$ws.Range("J7").Value = 0.03922895479591048
$ws.Range("O7").Value = 0.136569704933115
$ws.Range("P7").Value = 0.136569704933115
$ws.Range("S7").Value = 0.005357486781312003
$ws.Range("T7").Value = 0.005357486781312001

$ws.Range("I8").Value = 0.4989999962806363
$ws.Range("J8").Value = 0.4989999962806364
$ws.Range("M8").Value = 123.2806423333333
$ws.Range("N8").Value = 369.841927
$ws.Range("O8").Value = 0.6241574062367528
$ws.Range("P8").Value = 0.6241574062367526
$ws.Range("Q8").Value = 893.142077115278
$ws.Range("R8").Value = 8038.278694037502
$ws.Range("S8").Value = 0.3114545433906712
$ws.Range("T8").Value = 0.3114545433906712

$ws.Range("I9").Value = 0.4989999962806363
$ws.Range("J9").Value = 0.4989999962806364
$ws.Range("O9").Value = 0.2392728888301323
$ws.Range("P9").Value = 0.2392728888301322
$ws.Range("S9").Value = 0.1193971706362931
$ws.Range("T9").Value = 0.1193971706362931

$ws.Range("I10").Value = 0.4989999962806363
$ws.Range("J10").Value = 0.4989999962806364
$ws.Range("O10").Value = 0.136569704933115
$ws.Range("P10").Value = 0.136569704933115
$ws.Range("S10").Value = 0.068148282253672
$ws.Range("T10").Value = 0.068148282253672
